# Refresh the cryptocurrency Price (D) and Volume(1h) (E) columns with the
# latest scraped figures, row by row (rows 2-51), leaving the Coin/Link/
# rank columns untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column holds plain-text figures such as "1.003" or
# "28.153.39" (thousands-separated prices render with dots, not commas).
# Force it to Text format first so Excel does not reinterpret them as
# numbers/dates when we assign the new strings.
$ws.Range("D2:D51").NumberFormat = "@"

$data = New-Object "object[,]" 50,2
$data[0,0] = '28.153.39'
$data[0,1] = '  -1.31%  '
$data[1,0] = '1.794.00'
$data[1,1] = '  -1.60%  '
$data[2,0] = '1.003'
$data[2,1] = '  +0.24%  '
$data[3,0] = '314.25'
$data[3,1] = '  -0.27%  '
$data[4,0] = '1.002'
$data[4,1] = '  +0.18%  '
$data[5,0] = '0.5215'
$data[5,1] = '  +2.14%  '
$data[6,0] = '0.3812'
$data[6,1] = '  -3.42%  '
$data[7,0] = '0.07964'
$data[7,1] = '  -2.56%  '
$data[8,0] = '41.44'
$data[8,1] = '  -0.51%  '
$data[9,0] = '1.094'
$data[9,1] = '  -1.40%  '
$data[10,0] = '6.286'
$data[10,1] = '  -0.88%  '
$data[11,0] = '1.003'
$data[11,1] = '  +0.20%  '
$data[12,0] = '20.49'
$data[12,1] = '  -2.91%  '
$data[13,0] = '7.286'
$data[13,1] = '  -3.02%  '
$data[14,0] = '1.795.35'
$data[14,1] = '  -1.37%  '
$data[15,0] = '91.80'
$data[15,1] = '  -0.74%  '
$data[16,0] = '0.00001089'
$data[16,1] = '  -3.71%  '
$data[17,0] = '0.06561'
$data[17,1] = '  -1.49%  '
$data[18,0] = '1.002'
$data[18,1] = '  +0.22%  '
$data[19,0] = '17.31'
$data[19,1] = '  -2.88%  '
$data[20,0] = '5.949'
$data[20,1] = '  -2.37%  '
$data[21,0] = '28.197.72'
$data[21,1] = '  -1.22%  '
$data[22,0] = '11.13'
$data[22,1] = '  -2.20%  '
$data[23,0] = '2.269'
$data[23,1] = '  +0.22%  '
$data[24,0] = '160.36'
$data[24,1] = '  +2.68%  '
$data[25,0] = '20.43'
$data[25,1] = '  -4.23%  '
$data[26,0] = '1.997.75'
$data[26,1] = '  -1.56%  '
$data[27,0] = '2.336'
$data[27,1] = '  -2.79%  '
$data[28,0] = '122.60'
$data[28,1] = '  -2.61%  '
$data[29,0] = '0.1074'
$data[29,1] = '  -1.89%  '
$data[30,0] = '1.052'
$data[30,1] = '  -5.72%  '
$data[31,0] = '3.672'
$data[31,1] = '  +0.44%  '
$data[32,0] = '5.537'
$data[32,1] = '  -3.97%  '
$data[33,0] = '0.07229'
$data[33,1] = '  +2.30%  '
$data[34,0] = '12.09'
$data[34,1] = '  +6.99%  '
$data[35,0] = '0.02308'
$data[35,1] = '  -1.97%  '
$data[36,0] = '0.2142'
$data[36,1] = '  -3.80%  '
$data[37,0] = '8.708'
$data[37,1] = '  -1.52%  '
$data[38,0] = '5.064'
$data[38,1] = '  -3.54%  '
$data[39,0] = '0.6153'
$data[39,1] = '  -2.58%  '
$data[40,0] = '1.161'
$data[40,1] = '  -1.63%  '
$data[41,0] = '1.355'
$data[41,1] = '  -3.08%  '
$data[42,0] = '13.28'
$data[42,1] = '  -1.61%  '
$data[43,0] = '3.767'
$data[43,1] = '  +0.86%  '
$data[44,0] = '0.5946'
$data[44,1] = '  +0.35%  '
$data[45,0] = '127.84'
$data[45,1] = '  +2.16%  '
$data[46,0] = '1.225'
$data[46,1] = '  +3.33%  '
$data[47,0] = '1.916'
$data[47,1] = '  -3.47%  '
$data[48,0] = '0.06746'
$data[48,1] = '  -2.13%  '
$data[49,0] = '72.89'
$data[49,1] = '  -1.91%  '

$ws.Range("D2:E51").Value = $data

